$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85 (shifts old rows 85-134 down to 86-135)
$ws.Rows.Item(85).Insert()

# Fill in the new row 85 with the new weekly record
$ws.Cells.Item(85, 1).Value = 6
$ws.Cells.Item(85, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(85, 3).Value = "Metropolitana"
$ws.Cells.Item(85, 4).Value = 44873
$ws.Cells.Item(85, 5).Value = 13
$ws.Cells.Item(85, 6).Value = "Fruta"
$ws.Cells.Item(85, 7).Value = 100103
$ws.Cells.Item(85, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(85, 9).Value = 100103003
$ws.Cells.Item(85, 10).Value = "Damasco"
$ws.Cells.Item(85, 11).Value = "Castle Brite"
$ws.Cells.Item(85, 12).Value = "Primera"
$ws.Cells.Item(85, 13).Value = 120
$ws.Cells.Item(85, 14).Value = 28500
$ws.Cells.Item(85, 15).Value = 28500
$ws.Cells.Item(85, 16).Value = 28500
$ws.Cells.Item(85, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(85, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(85, 19).Value = 2850
$ws.Cells.Item(85, 20).Value = 10
